$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" / "Valor Mora" table (rows 16-56) was listed newest
# period first (1911) down to oldest (1607). Re-sort it chronologically
# ascending (oldest 1607 first ... newest 1911 last) so new statement
# periods can be appended afterwards, while keeping each row's existing
# cell formatting (borders, shading, etc.) fixed in place.
$firstRow = 16
$lastRow = 56

$periods = @()
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $values += $ws.Cells.Item($r, 6).Value2
}

$revPeriods = @()
$revValues = @()
for ($i = $periods.Count - 1; $i -ge 0; $i--) {
    $revPeriods += $periods[$i]
    $revValues += $values[$i]
}

$i = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $revPeriods[$i]
    $ws.Cells.Item($r, 6).Value = $revValues[$i]
    $i++
}
